# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.836.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.23%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.117.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.02%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.54%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.111.57"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.13%  "

# Row 9
$ws.Range("E9").Value = "  -0.38%  "

# Row 10
$ws.Range("E10").Value = "  +11.09%  "

# Row 11
$ws.Range("E11").Value = "  +0.18%  "

# Row 12
$ws.Range("E12").Value = "  -0.60%  "

# Row 13
$ws.Range("E13").Value = "  +3.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.02%  "

# Row 15
$ws.Range("E15").Value = "  -1.05%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.633.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.701.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.11%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.01%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.118.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.45%  "

# Row 21
$ws.Range("E21").Value = "  +1.48%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.57%  "

# Row 23
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("E24").Value = "  -3.74%  "

# Row 25
$ws.Range("E25").Value = "  +0.14%  "

# Row 26
$ws.Range("E26").Value = "  -0.07%  "

# Row 27
$ws.Range("E27").Value = "  +7.97%  "

# Row 28
$ws.Range("E28").Value = "  -0.26%  "

# Row 29
$ws.Range("E29").Value = "  -1.42%  "

# Row 30
$ws.Range("E30").Value = "  -0.11%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "

# Row 33
$ws.Range("E33").Value = "  -2.87%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0876"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.87%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.50%  "

# Row 36
$ws.Range("E36").Value = "  +0.92%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.80%  "

# Row 38
$ws.Range("E38").Value = "  +0.02%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "448.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.07%  "

# Row 41
$ws.Range("E41").Value = "  -1.10%  "

# Row 42
$ws.Range("E42").Value = "  -1.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.886.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.19%  "

# Row 44
$ws.Range("E44").Value = "  -0.44%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.26%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.49%  "

# Row 48
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "

# Row 49
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.40%  "

# Row 50
$ws.Range("E50").Value = "  -0.54%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.68%  "
